# ------------------------------------------------------------------
# Data_Catalog.xlsx reorganization
#   - insert "NumFile" column after "Dataset"
#   - insert "File(s)" column before "Description"
#   - add "Reference" column at the end
#   - reword a few Description cells
#   - add a new "2025-Metabotypes" dataset spanning two rows (two files)
#   - add a hyperlink in the Reference column for the first new row
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# alignment constants (late bound COM, avoid relying on enum types)
$xlCenter = -4108
$xlLeft   = -4131

# ---------------------------------------------------------------
# 1. Header row
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Dataset"
$ws.Range("B1").Value = "NumFile"
$ws.Range("C1").Value = "Samples"
$ws.Range("D1").Value = "Features"
$ws.Range("E1").Value = "File(s)"
$ws.Range("F1").Value = "Description"
$ws.Range("G1").Value = "Reference"

$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").VerticalAlignment = $xlCenter

$ws.Range("B1:E1").HorizontalAlignment = $xlCenter
$ws.Range("F1").WrapText = $true

$ws.Range("G1").Font.Bold = $false
$ws.Range("G1").VerticalAlignment = $xlCenter
$ws.Range("G1").WrapText = $true

# ---------------------------------------------------------------
# 2. Data rows (A2:G7 keep the same datasets, new columns filled in)
# ---------------------------------------------------------------

# Row 2 - 2018-MetabotypingPaper
$ws.Range("A2").Value = "2018-MetabotypingPaper"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 39
$ws.Range("D2").Value = 690
$ws.Range("E2").Value = "DataValues_S013.csv"
$ws.Range("F2").Value = 'Data used in the paper "Metabotypes of response to bariatric surgery independent of the magnitude of weight loss"'
$ws.Range("G2").Value = "https://doi.org/10.1371/journal.pone.0198214"
$ws.Hyperlinks.Add($ws.Range("G2"), "https://doi.org/10.1371/journal.pone.0198214") | Out-Null

# Row 3 - 2018-Phosphoproteomics
$ws.Range("A3").Value = "2018-Phosphoproteomics"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 1320
$ws.Range("E3").Value = "TIO2+PTYR-human-MSS+MSIvsPD.XLSX"
$ws.Range("F3").Value = "Data obtained from a phosphoproteomics experiment that was performed to analyze (3 + 3) PDX models of two different subtypes using Phosphopeptide enriched samples."

# Row 4 - 2023-CIMCBTutorial
$ws.Range("A4").Value = "2023-CIMCBTutorial"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 140
$ws.Range("D4").Value = 149
$ws.Range("E4").Value = "GastricCancer_NMR.xlsx"
$ws.Range("F4").Value = 'NMR data from a gastric cancer study used in a metabolomics data analysis tutorial ("Basic Metabolomics Data Analysis Workflow" (https://cimcb.github.io/MetabWorkflowTutorial/Tutorial1.html)'

# Row 5 - 2023-UGrX-4MetaboAnalystTutorial
$ws.Range("A5").Value = "2023-UGrX-4MetaboAnalystTutorial"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 24
$ws.Range("D5").Value = 145
$ws.Range("E5").Value = "ST000002_AN000002_clean.csv"
$ws.Range("F5").Value = "Data from MetabolomicsWorkbench (ID ST000002)"

# Row 6 - 2024-fobitools-UseCase_1
$ws.Range("A6").Value = "2024-fobitools-UseCase_1"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 45
$ws.Range("D6").Value = 1541
$ws.Range("E6").Value = "ST000291curated.xlsx"
$ws.Range("F6").Value = 'Data used in the fobitools Bioconductor package, in one its vignettes, [Use Case ST000291] analyzing the data from Metabolomics Workbench  Dataset '

# Row 7 - 2024-Cachexia
$ws.Range("A7").Value = "2024-Cachexia"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 77
$ws.Range("D7").Value = 63
$ws.Range("E7").Value = "human_cachexia.csv"
$ws.Range("F7").Value = 'Data used in several MetaboAnalyst tutorials. 77 urine samples, 47 patients with cachexia, and 30 control patients (from the "specmine.datasets" R package)'

# Row 8 - 2025-Metabotypes (file 1 of 2)
$ws.Range("A8").Value = "2025-Metabotypes"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 110
$ws.Range("D8").Value = 892
$ws.Range("E8").Value = "data/results_filtered_normalized.csv"
$ws.Range("F8").Value = 'filtered and normalized data as described in the acompanying document "code/Reanálisis_de_resultados_procesados-Estudio_ST002993.html"'
$ws.Range("G8").Value = '(probably part of) Data used in the paper "Identifying subgroups of childhood obesity by using multiplatform metabotyping"'

# Row 9 - 2025-Metabotypes (file 2 of 2)
$ws.Range("A9").Value = "2025-Metabotypes"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 110
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = "data/factor_matrix.csv"
$ws.Range("F9").Value = 'Dataset with 5 (latent factors) x 110 (samples)   matrix produced by the multiple factor analysis as described in the acompanying document "code/Reanálisis_de_resultados_procesados-Estudio_ST002993.html"'

# ---------------------------------------------------------------
# 3. Formatting that mirrors the original workbook's conventions
# ---------------------------------------------------------------

# Column B / C / D data cells: vertical-center (matches the sheet's
# existing convention for the numeric columns)
$ws.Range("B2:D9").VerticalAlignment = $xlCenter

# Column E (File(s)) - left aligned, vertical centered
$ws.Range("E2:E9").HorizontalAlignment = $xlLeft
$ws.Range("E2:E9").VerticalAlignment = $xlCenter

# Column F (Description) and G (Reference) - wrap text, vertical centered
$ws.Range("F2:G9").WrapText = $true
$ws.Range("F2:G9").VerticalAlignment = $xlCenter

# Row heights (wrap-text driven row heights in the source workbook)
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45

# Column widths (values chosen so the resulting serialized width lands on
# the target character width as closely as this engine's rounding allows)
# Column A keeps its original width/bestFit - left untouched on purpose.
$ws.Columns.Item(2).ColumnWidth = 8.5
$ws.Columns.Item(3).ColumnWidth = 10.666666666666
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 34.666666666666
$ws.Columns.Item(6).ColumnWidth = 76.0
$ws.Columns.Item(7).ColumnWidth = 41.166666666666

# ---------------------------------------------------------------
# 4. View state
# ---------------------------------------------------------------
$ws.Range("F10").Select()
